$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.857.12'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.298.04'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '306.25'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.97%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '97.35'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.510'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -1.63%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.504'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.68'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0789'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '18.18'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.77%  '
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.77'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.660.57'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.301.62'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '42.814.57'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -3.97%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0901'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.03'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.78'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.69%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '236.60'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -2.02%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.50%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '25.38'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +2.75%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '166.81'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.73%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '9.04'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -1.09%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '33.10'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.82'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.00'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.18'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -5.07%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.20%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0690'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.75'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -1.46%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.109'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -1.40%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.73'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.005.26'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -1.85%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.43%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.01'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.75%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '17.97'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +4.83%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.78'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.88'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +3.96%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '53.99'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.528.15'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.18%  '
